$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell A1 changes from "Witt and Jackson 2016" to "name"
$ws.Range("A1").Value = "name"

# Add a new row 24 with "hahaha" in column A
$ws.Range("A24").Value = "hahaha"

# Update the selection / active cell to D14 (as recorded in the saved view)
$ws.Range("D14").Select()
